$d = $word.ActiveDocument

function Accept-AllPendingRevisions($doc) {
    # Accept revisions one at a time (back to front) instead of calling
    # Document.AcceptAllRevisions(), which also strips the w:rsid* trail
    # from every run in the package as a side effect. Accepting the
    # individual Revision objects only resolves the edits we just made.
    for ($i = $doc.Revisions.Count; $i -ge 1; $i--) {
        $doc.Revisions.Item($i).Accept()
    }
}

# Track changes + per-revision Accept keeps the inserted text in its own
# <w:r> instead of silently re-merging it into the neighbouring run that
# already carries identical formatting (this mirrors what a real editing
# session in Word produces).
$d.TrackRevisions = $true

# --- 1. "April 2017" -> "June 12, " + "2017" (kept as two runs) ---------
$rng1 = $d.Content
$rng1.Find.Execute("April", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$aprilRange = $rng1.Duplicate
$dateIns = $d.Range($aprilRange.Start, $aprilRange.Start)
$dateIns.InsertBefore("June 12, ")

$rng2 = $d.Content
$rng2.Find.Execute("April ", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$oldDate = $rng2.Duplicate
$oldDate.Delete()

Accept-AllPendingRevisions $d

# --- 2. Relocate the _GoBack bookmark from the "We suggest " paragraph --
#        to the middle of "experimental" (split "exp" | "erimental ..."),
#        matching where the author's cursor ended up after their edit.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$rng3 = $d.Content
$rng3.Find.Execute("inspire innovative exp", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$splitHit = $rng3.Duplicate
$splitPos = $splitHit.End

# Force a run break exactly at splitPos without changing any visible text.
$tick = $d.Range($splitPos, $splitPos)
$tick.InsertAfter("X")
$tickSel = $d.Range($splitPos, $splitPos + 1)
$tickSel.Delete()

Accept-AllPendingRevisions $d
$d.TrackRevisions = $false

$newGoBackRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $newGoBackRange) | Out-Null
